$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row at row 89 (pushes existing rows 89..186 down to 90..187)
$ws.Rows(89).Insert()

$ws.Range("A89").Value = 4
$ws.Range("B89").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C89").Value = "Los Lagos"
$ws.Range("D89").Value = 44539
$ws.Range("E89").Value = 10
$ws.Range("F89").Value = 100112044
$ws.Range("G89").Value = "Perejil"
$ws.Range("H89").Value = "Sin especificar"
$ws.Range("I89").Value = "Primera"
$ws.Range("J89").Value = 60
$ws.Range("K89").Value = 5000
$ws.Range("L89").Value = 5000
$ws.Range("M89").Value = 5000
$ws.Range("N89").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O89").Value = "Región Metropolitana"
$ws.Range("P89").Value = 1667
$ws.Range("Q89").Value = 3
$ws.Range("R89").Value = "Hortaliza"
